$d = $word.ActiveDocument

# 1. Replace the placeholder ID text and drop the trailing single-space run.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5301_topic_48__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SMC_PGI_5301_91__ID**", 2)

# 2. Update the first paragraph's formatting: indent + paragraph border.
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
